$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -6
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -2
$ws.Range("F9").Value = -4
$ws.Range("F13").Value = 2
